# Apply "want to go" / price updates to each sheet of the workbook, as
# generated by the site's data refresh (gh-pages output regenerated at
# 456a3b4). Sheet order in the workbook is:
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Performances)
#   3 = 本地生活 (Local life)
#   4 = 全部类型 (All types - union of the above three)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 866
$ws1.Range("F4").Value  = 1092
$ws1.Range("F5").Value  = 518
$ws1.Range("F6").Value  = 219
$ws1.Range("F7").Value  = 663
$ws1.Range("F8").Value  = 243
$ws1.Range("F10").Value = 84
$ws1.Range("F11").Value = 216
$ws1.Range("F12").Value = 149
$ws1.Range("F13").Value = 1814
$ws1.Range("G13").Value = 88
$ws1.Range("F14").Value = 428
$ws1.Range("F15").Value = 41
$ws1.Range("F16").Value = 490
$ws1.Range("F17").Value = 257
$ws1.Range("F21").Value = 661
$ws1.Range("F22").Value = 47
$ws1.Range("F23").Value = 242
$ws1.Range("F24").Value = 956
$ws1.Range("F26").Value = 1549
$ws1.Range("F27").Value = 295

# --- Sheet 2: 演出 -------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 11
$ws2.Range("F8").Value = 280

# --- Sheet 3: 本地生活 ---------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 394

# --- Sheet 4: 全部类型 (aggregated view) ---------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 394
$ws4.Range("F3").Value  = 866
$ws4.Range("F5").Value  = 1092
$ws4.Range("F8").Value  = 518
$ws4.Range("F9").Value  = 219
$ws4.Range("F10").Value = 663
$ws4.Range("F12").Value = 243
$ws4.Range("F14").Value = 84
$ws4.Range("F15").Value = 216
$ws4.Range("F16").Value = 149
$ws4.Range("F17").Value = 1814
$ws4.Range("G17").Value = 88
$ws4.Range("F19").Value = 428
$ws4.Range("F20").Value = 41
$ws4.Range("F21").Value = 490
$ws4.Range("F22").Value = 257
$ws4.Range("F27").Value = 11
$ws4.Range("F28").Value = 280
$ws4.Range("F30").Value = 661
$ws4.Range("F35").Value = 47
$ws4.Range("F36").Value = 242
$ws4.Range("F37").Value = 956
$ws4.Range("F39").Value = 1549
$ws4.Range("F40").Value = 295
